# Update "Activity Relationship table.xlsx"
# - Rename the "Outside" activity to "Atrium" (row 14 label + matching column header N1)
# - Update several relationship scores that changed alongside that rename

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Outside" -> "Atrium" (row label in column A and matching column header in row 1)
$ws.Range("A14").Value = "Atrium"
$ws.Range("N1").Value = "Atrium"

# Updated relationship values against the renamed "Atrium" activity (column N)
$ws.Range("N3").Value = 4
$ws.Range("N4").Value = 4
$ws.Range("N6").Value = 4
$ws.Range("N7").Value = 4

# Updated relationship values on the "Atrium" row (row 14)
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = 4
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 4

# Update the active selection to match the saved state
$ws.Range("J22").Select() | Out-Null
